$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "D2" "66.833.13"
$ws.Range("E2").Value = "  -0.05%  "
Set-TextValue "D3" "3.784.79"
$ws.Range("E3").Value = "  -1.83%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue "D5" "439.03"
$ws.Range("E5").Value = "  +3.61%  "
Set-TextValue "D6" "141.34"
$ws.Range("E6").Value = "  +8.48%  "
Set-TextValue "D7" "0.623"
$ws.Range("E7").Value = "  +2.41%  "
$ws.Range("E8").Value = "  +0.00%  "
Set-TextValue "D9" "0.736"
$ws.Range("E9").Value = "  +1.59%  "
Set-TextValue "D10" "0.153"
$ws.Range("E10").Value = "  -7.64%  "
Set-TextValue "D11" "0.0000317"
$ws.Range("E11").Value = "  -10.44%  "
Set-TextValue "D12" "42.94"
$ws.Range("E12").Value = "  +5.86%  "
Set-TextValue "D13" "10.46"
$ws.Range("E13").Value = "  +4.81%  "
Set-TextValue "D14" "4.392.73"
$ws.Range("E14").Value = "  -1.91%  "
Set-TextValue "D15" "14.80"
$ws.Range("E15").Value = "  -9.98%  "
$ws.Range("E16").Value = "  -0.66%  "
Set-TextValue "D17" "3.833.60"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("E19").Value = "  +6.53%  "
Set-TextValue "D20" "66.936.55"
$ws.Range("E20").Value = "  -0.20%  "
Set-TextValue "D21" "413.38"
$ws.Range("E21").Value = "  +2.08%  "
Set-TextValue "D22" "14.49"
$ws.Range("E22").Value = "  +1.10%  "
Set-TextValue "D23" "3.26"
$ws.Range("E23").Value = "  +8.27%  "
Set-TextValue "D24" "85.71"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D25" "3.43"
$ws.Range("E25").Value = "  +7.75%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D26" "37.05"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("E27").Value = "  -0.78%  "
$ws.Range("E28").Value = "  +3.13%  "
Set-TextValue "D29" "9.51"
$ws.Range("E29").Value = "  +31.67%  "
Set-TextValue "D30" "731.24"
$ws.Range("E30").Value = "  +4.77%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D31" "13.93"
$ws.Range("E31").Value = "  +12.97%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D32" "0.133"
$ws.Range("E32").Value = "  +10.10%  "
Set-TextValue "D33" "2.73"
$ws.Range("E33").Value = "  -1.12%  "
Set-TextValue "D34" "42.78"
$ws.Range("E34").Value = "  +11.64%  "
Set-TextValue "D35" "0.157"
$ws.Range("E35").Value = "  +3.70%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D36" "5.63"
$ws.Range("E36").Value = "  +27.50%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D37" "56.46"
$ws.Range("E37").Value = "  +2.36%  "
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D38" "1.00"
$ws.Range("E38").Value = "  +0.12%  "
Set-TextValue "D39" "0.0477"
$ws.Range("E39").Value = "  +4.67%  "
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D40" "2.92"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D41" "2.62"
$ws.Range("E41").Value = "  +32.40%  "
Set-TextValue "D42" "0.0₃0677"
$ws.Range("E42").Value = "  -14.04%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D43" "1.00"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D44" "0.140"
$ws.Range("E44").Value = "  +3.54%  "
$ws.Range("E45").Value = "  +6.54%  "
Set-TextValue "D46" "0.327"
$ws.Range("E46").Value = "  +13.57%  "
$ws.Range("E47").Value = "  +1.20%  "
$ws.Range("E48").Value = "  +4.12%  "
Set-TextValue "D49" "2.09"
$ws.Range("E49").Value = "  +1.41%  "
Set-TextValue "D50" "142.64"
$ws.Range("E50").Value = "  -2.57%  "
Set-TextValue "D51" "2.82"
$ws.Range("E51").Value = "  +1.80%  "
